$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B49: change from text "1" to a real numeric value 1
$ws.Range("B49").Value = 1

# Add new row 50 with annotation data
$ws.Range("A50").Value = "Ying Tang"

# B50 must remain textual ("3"), not get auto-converted to a number.
# Temporarily format the cell as Text so Excel keeps it as a string,
# then reset the cell style back to Normal so no stray style is left
# on the cell while the value stays stored as text.
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "3"
$ws.Range("B50").Style = "Normal"

$ws.Range("C50").Value = "无"
$ws.Range("D50").Value = "DFT"
$ws.Range("E50").Value = "WRI"
$ws.Range("F50").Value = "3bf3a8cd-f7a3-492e-815a-c1d9e74634b1"
$ws.Range("G50").Value = "ByCPHrgCW_annotated.xlsx"
$ws.Range("H50").Value = "The problem scenario states that the model/weights is private, but later on it ceases to be so (weights are not encrypted)."
